$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the environment host (dev4 -> dev4-oci).
# A2 keeps its original "quote prefix" (text stored via a leading apostrophe)
# cell style, so re-enter it the same way to avoid losing that formatting.
$ws.Range("A2").Value = "'ssurgwsoadev4-oci.opc.oracleoutsourcing.com"

# Update the URL cell + its hyperlink target (dev4 -> dev4-oci).
# The engine's Hyperlinks.Add always mints a fresh relationship rather than
# patching one in place, so drop the stale link first to avoid a duplicate.
$newUrl = "https://ssurgwsoadev4-oci.opc.oracleoutsourcing.com/pc/PolicyCenter.do"
$ws.Range("B2").Hyperlinks.Delete()
$ws.Range("B2").Value = $newUrl
$ws.Hyperlinks.Add($ws.Range("B2"), $newUrl) | Out-Null

# Update PAS number
$ws.Range("G2").Value = 6965

# Update NombreProductor
$ws.Range("H2").Value = "Corporativos Directos"

# Update the active selection to F2, as the last-saved view
$ws.Range("F2").Select()
